$d = $word.ActiveDocument

# 1. Heading2 title: "Some Class Name" -> "ForthInterpreter"
$d.Content.Find.Execute("Some Class Name", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ForthInterpreter", 2)

# 2. Table cell: "Class Name" -> "ForthInterpreter"
$d.Content.Find.Execute("Class Name", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ForthInterpreter", 2)

# 3. Merge "So" + bookmark + "me paragraph about this class" into
#    "Some paragraph about this class" (bookmark removed from this location)
$d.Content.Find.Execute("Some paragraph about this class Some paragraph about this class", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Some paragraph about this class Some paragraph about this class", 2)
